# Applies the "Updated cryptos list" data refresh described by the diff:
#   - Column D (Price) and Column E (Volume(1h)) values updated for most rows.
#   - Three coin pairs swapped rows (name/link/price/volume moved as a unit):
#       row 16 <-> row 17   (Avalanche  <-> ShibaInu)
#       row 28 <-> row 29   (WrappedeETH <-> Kaspa)
#       row 43 <-> row 44   (Hedera <-> Filecoin)
#
# All D/E cells in this sheet are stored as TEXT (t="inlineStr"), even when the
# text looks like a plain number (e.g. "5.38"). A bare Range.Value assignment of
# such a string gets auto-coerced to a numeric cell by Excel, which would change
# the cell type and not match the source. Set-TextValue forces text entry (via a
# temporary "@" number format) and then resets the style back to Normal so no
# stray style/format is left behind on the cell.

function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.412.50'
$ws.Range('E2').Value = '  -8.03%  '
$ws.Range('D3').Value = '2.437.33'
$ws.Range('E3').Value = '  -13.77%  '
Set-TextValue $ws 'D5' '466.65'
$ws.Range('E5').Value = '  -7.37%  '
Set-TextValue $ws 'D6' '131.96'
$ws.Range('E6').Value = '  -2.72%  '
Set-TextValue $ws 'D7' '0.994'
$ws.Range('E7').Value = '  -0.56%  '
Set-TextValue $ws 'D8' '0.492'
$ws.Range('E8').Value = '  -7.42%  '
$ws.Range('D9').Value = '2.449.87'
$ws.Range('E9').Value = '  -13.20%  '
Set-TextValue $ws 'D10' '0.0963'
$ws.Range('E10').Value = '  -6.82%  '
Set-TextValue $ws 'D11' '5.38'
$ws.Range('E11').Value = '  -9.58%  '
Set-TextValue $ws 'D12' '0.320'
$ws.Range('E12').Value = '  -8.17%  '
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('D14').Value = '2.856.43'
$ws.Range('E14').Value = '  -14.01%  '
$ws.Range('D15').Value = '54.382.95'
$ws.Range('E15').Value = '  -8.21%  '
# Row 16/17: Avalanche and ShibaInu swap places
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 'D16' '0.0000134'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws 'D17' '19.85'
$ws.Range('E17').Value = '  -8.49%  '

$ws.Range('D18').Value = '2.451.94'
$ws.Range('E18').Value = '  -12.09%  '
$ws.Range('E19').Value = '  -10.73%  '
Set-TextValue $ws 'D20' '314.11'
$ws.Range('E20').Value = '  -10.83%  '
Set-TextValue $ws 'D21' '9.41'
$ws.Range('E21').Value = '  -15.03%  '
$ws.Range('E22').Value = '  +0.49%  '
Set-TextValue $ws 'D23' '5.68'
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('E24').Value = '  -13.72%  '
Set-TextValue $ws 'D25' '56.91'
$ws.Range('E25').Value = '  -9.90%  '
$ws.Range('E26').Value = '  -3.30%  '
Set-TextValue $ws 'D27' '0.385'
$ws.Range('E27').Value = '  -10.35%  '
# Row 28/29: WrappedeETH and Kaspa swap places
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D28' '0.154'
$ws.Range('E28').Value = '  -10.81%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.515.91'
$ws.Range('E29').Value = '  -14.90%  '

Set-TextValue $ws 'D30' '7.20'
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').Value = '0.0₃0737'
$ws.Range('E32').Value = '  -8.87%  '
Set-TextValue $ws 'D33' '151.02'
$ws.Range('E33').Value = '  +0.47%  '
Set-TextValue $ws 'D34' '17.70'
$ws.Range('E35').Value = '  -11.53%  '
Set-TextValue $ws 'D36' '5.07'
$ws.Range('E36').Value = '  -5.32%  '
Set-TextValue $ws 'D37' '3.55'
$ws.Range('E37').Value = '  -14.84%  '
Set-TextValue $ws 'D38' '1.06'
$ws.Range('E38').Value = '  -6.36%  '
Set-TextValue $ws 'D39' '0.802'
$ws.Range('E39').Value = '  -10.90%  '
Set-TextValue $ws 'D40' '33.61'
$ws.Range('E40').Value = '  -8.15%  '
$ws.Range('E41').Value = '  -0.58%  '
Set-TextValue $ws 'D42' '0.608'
$ws.Range('E42').Value = '  -2.91%  '
# Row 43/44: Hedera and Filecoin swap places
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D43' '3.30'
$ws.Range('E43').Value = '  -7.06%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D44' '0.0530'
$ws.Range('E44').Value = '  -5.40%  '

Set-TextValue $ws 'D45' '10.20'
$ws.Range('E45').Value = '  -1.39%  '
Set-TextValue $ws 'D46' '1.25'
$ws.Range('E46').Value = '  -8.34%  '
$ws.Range('D47').Value = '1.976.02'
$ws.Range('E47').Value = '  -10.58%  '
Set-TextValue $ws 'D48' '0.0221'
$ws.Range('E48').Value = '  -2.03%  '
Set-TextValue $ws 'D49' '0.0871'
$ws.Range('E49').Value = '  -1.89%  '
Set-TextValue $ws 'D50' '4.42'
$ws.Range('E50').Value = '  -2.37%  '
Set-TextValue $ws 'D51' '16.66'
$ws.Range('E51').Value = '  -13.69%  '
